# repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to reflect the
# re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -6
    4  = -8
    5  = -12
    9  = -4
    12 = 2
    18 = -6
    22 = 1
    23 = 0
    27 = -5
    29 = -3
    30 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
